$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.536.36"
$ws.Range("D3").Value = "2.510.45"
$ws.Range("E3").Value = "  -4.95%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.71%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "2.509.07"
$ws.Range("E9").Value = "  -5.01%  "
$ws.Range("E10").Value = "  -7.59%  "
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "2.967.92"
$ws.Range("E14").Value = "  -5.03%  "
$ws.Range("D15").Value = "69.390.33"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("E16").Value = "  -6.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").Value = "2.502.49"
$ws.Range("E18").Value = "  -6.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.71%  "
$ws.Range("E22").Value = "  -3.96%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  -5.97%  "
$ws.Range("E27").Value = "  -5.86%  "
$ws.Range("D28").Value = "2.637.45"
$ws.Range("E28").Value = "  -5.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.993"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "463.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.34%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.119"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.79%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.04%  "
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("E43").Value = "  -6.76%  "
$ws.Range("E44").Value = "  -14.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.42%  "
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("E49").Value = "  -4.12%  "
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("E51").Value = "  -1.93%  "
